$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 3571.0527
$ws.Range("I5").Value = 1006.4545
$ws.Range("J5").Value = 7097.375
$ws.Range("K5").Value = 1006.4545
$ws.Range("L5").Value = 7097.375
$ws.Range("M5").Value = -891.4545
$ws.Range("N5").Value = -7327.375

$ws.Range("H103").Value = 1543.4783
$ws.Range("J103").Value = 1727.8182
$ws.Range("L103").Value = 5183.4546
$ws.Range("N103").Value = -6355.4546

$ws.Range("H129").Value = 11906463
$ws.Range("I129").Value = 830.75
$ws.Range("J129").Value = 20835688
$ws.Range("K129").Value = 2492.25
$ws.Range("L129").Value = 62507064
$ws.Range("M129").Value = 2507.75
$ws.Range("N129").Value = -62517064

$ws.Range("H131").Value = 1062.3334
$ws.Range("I131").Value = 1062.3334
$ws.Range("K131").Value = 3187.0002
$ws.Range("M131").Value = 1852.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 454
$ws.Range("I4").Value = 412
$ws.Range("K4").Value = 412
$ws.Range("M4").Value = -296

$ws.Range("H61").Value = 2872.647
$ws.Range("J61").Value = 3582.5
$ws.Range("L61").Value = 3582.5
$ws.Range("N61").Value = -4006.5

$ws.Range("H63").Value = 2998.4119
$ws.Range("I63").Value = 1396.625
$ws.Range("K63").Value = 1396.625
$ws.Range("M63").Value = -710.625

$ws.Range("H66").Value = 2998.4119
$ws.Range("I66").Value = 1396.625
$ws.Range("K66").Value = 6983.125
$ws.Range("M66").Value = -3551.125

$ws.Range("H74").Value = 2069.2222
$ws.Range("I74").Value = 2224.8
$ws.Range("J74").Value = 1874.75
$ws.Range("K74").Value = 2224.8
$ws.Range("L74").Value = 1874.75
$ws.Range("M74").Value = -1350.8
$ws.Range("N74").Value = -3622.75

$ws.Range("H77").Value = 2069.2222
$ws.Range("I77").Value = 2224.8
$ws.Range("J77").Value = 1874.75
$ws.Range("K77").Value = 11124
$ws.Range("L77").Value = 9373.75
$ws.Range("M77").Value = -6756
$ws.Range("N77").Value = -18109.75

$ws.Range("H98").Value = 45000
$ws.Range("J98").Value = 45000
$ws.Range("L98").Value = 45000
$ws.Range("N98").Value = -50990

$ws.Range("H101").Value = 32499
$ws.Range("J101").Value = 32499
$ws.Range("L101").Value = 32499
$ws.Range("N101").Value = -38989

$ws.Range("H132").Value = 2605683
$ws.Range("I132").Value = 3790587.5
$ws.Range("J132").Value = 235874.5
$ws.Range("K132").Value = 11371762.5
$ws.Range("L132").Value = 707623.5
$ws.Range("M132").Value = -11369232.5
$ws.Range("N132").Value = -712683.5

$ws.Range("H136").Value = 2872.647
$ws.Range("J136").Value = 3582.5
$ws.Range("L136").Value = 10747.5
$ws.Range("N136").Value = -15847.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 4605.5713
$ws.Range("I64").Value = 1698
$ws.Range("J64").Value = 11874.5
$ws.Range("K64").Value = 1698
$ws.Range("L64").Value = 11874.5
$ws.Range("M64").Value = -1473
$ws.Range("N64").Value = -12324.5

$ws.Range("H67").Value = 4605.5713
$ws.Range("I67").Value = 1698
$ws.Range("J67").Value = 11874.5
$ws.Range("K67").Value = 1698
$ws.Range("L67").Value = 11874.5
$ws.Range("M67").Value = -918
$ws.Range("N67").Value = -13434.5

$ws.Range("H86").Value = 3039.08
$ws.Range("I86").Value = 2106.2
$ws.Range("K86").Value = 2106.2
$ws.Range("M86").Value = -983.1999999999998

$ws.Range("H89").Value = 3039.08
$ws.Range("I89").Value = 2106.2
$ws.Range("K89").Value = 10531
$ws.Range("M89").Value = -4915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 210.76923
$ws.Range("I7").Value = 220
$ws.Range("K7").Value = 220
$ws.Range("M7").Value = -107

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 4261.4287
$ws.Range("I122").Value = 3908
$ws.Range("K122").Value = 11724
$ws.Range("M122").Value = -9274

$ws.Range("H134").Value = 71434760
$ws.Range("I134").Value = 166668660
$ws.Range("J134").Value = 9336.875
$ws.Range("K134").Value = 500005980
$ws.Range("L134").Value = 28010.625
$ws.Range("M134").Value = -500003445
$ws.Range("N134").Value = -33080.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 143828.45
$ws.Range("I5").Value = 901.2222
$ws.Range("K5").Value = 2703.6666
$ws.Range("M5").Value = -2591.6666

$ws.Range("H68").Value = 386725.78
$ws.Range("J68").Value = 501673.6
$ws.Range("L68").Value = 1505020.8
$ws.Range("N68").Value = -1506642.8

$ws.Range("H71").Value = 386725.78
$ws.Range("J71").Value = 501673.6
$ws.Range("L71").Value = 4515062.399999999
$ws.Range("N71").Value = -4523174.399999999

$ws.Range("H87").Value = 20871.625
$ws.Range("I87").Value = 9493.25
$ws.Range("K87").Value = 28479.75
$ws.Range("M87").Value = -27231.75

$ws.Range("H90").Value = 20871.625
$ws.Range("I90").Value = 9493.25
$ws.Range("K90").Value = 85439.25
$ws.Range("M90").Value = -79199.25

$ws.Range("H92").Value = 198.9
$ws.Range("J92").Value = 127.75
$ws.Range("L92").Value = 383.25
$ws.Range("N92").Value = -2879.25

$ws.Range("H111").Value = 6618.6665
$ws.Range("J111").Value = 11537.25
$ws.Range("L111").Value = 34611.75
$ws.Range("N111").Value = -40745.75

$ws.Range("H121").Value = 119095
$ws.Range("I121").Value = 731.6667
$ws.Range("J121").Value = 163481.25
$ws.Range("K121").Value = 2195.0001
$ws.Range("L121").Value = 490443.75
$ws.Range("M121").Value = -885.0001000000002
$ws.Range("N121").Value = -493063.75

$ws.Range("H134").Value = 4127.6763
$ws.Range("J134").Value = 6999.1665
$ws.Range("L134").Value = 20997.4995
$ws.Range("N134").Value = -31137.4995

$ws.Range("H135").Value = 143828.45
$ws.Range("I135").Value = 901.2222
$ws.Range("K135").Value = 8110.999800000001
$ws.Range("M135").Value = -5575.999800000001

$ws.Range("H140").Value = 26044660
$ws.Range("I140").Value = 28738418
$ws.Range("K140").Value = 86215254
$ws.Range("M140").Value = -86210074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5579.684
$ws.Range("I7").Value = 3797.8
$ws.Range("K7").Value = 3797.8
$ws.Range("M7").Value = -3685.8

$ws.Range("H55").Value = 2320.8845
$ws.Range("I55").Value = 958
$ws.Range("J55").Value = 3683.7693
$ws.Range("K55").Value = 958
$ws.Range("L55").Value = 3683.7693
$ws.Range("M55").Value = -785
$ws.Range("N55").Value = -4029.7693

$ws.Range("H93").Value = 1055.7222
$ws.Range("I93").Value = 539.6
$ws.Range("K93").Value = 539.6
$ws.Range("M93").Value = 708.4

$ws.Range("H122").Value = 4351.2383
$ws.Range("I122").Value = 3854.5557
$ws.Range("J122").Value = 4723.75
$ws.Range("K122").Value = 11563.6671
$ws.Range("L122").Value = 14171.25
$ws.Range("M122").Value = -9113.667099999999
$ws.Range("N122").Value = -19071.25

$ws.Range("H126").Value = 5579.684
$ws.Range("I126").Value = 3797.8
$ws.Range("K126").Value = 11393.4
$ws.Range("M126").Value = -8923.400000000001
